# Update the "取得日時" (acquired datetime) column on the first sheet
# ("ランサーズ") for existing data rows to reflect the latest fetch
# timestamp, 2025-10-02 12:44:43 (append run at 2025-10-02 12:44 JST).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-02 12:44:43"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
